# Update NATMI Cd38-Pecam1 LR-pair stats with recomputed TPM-based values
# (ligand/receptor expression, specificity, and edge-weight columns G:T
# for data rows 2-17), per the "update scripts wuth new tpm" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = [double]"60.813934"
$ws.Cells.Item(2, 8).Value = [double]"182.441802"
$ws.Cells.Item(2, 9).Value = [double]"0.6840634102070431"
$ws.Cells.Item(2, 10).Value = [double]"0.6840634102070431"
$ws.Cells.Item(2, 13).Value = [double]"449.104309"
$ws.Cells.Item(2, 14).Value = [double]"1347.312927"
$ws.Cells.Item(2, 15).Value = [double]"0.9710020245482639"
$ws.Cells.Item(2, 16).Value = [double]"0.9710020245482639"
$ws.Cells.Item(2, 17).Value = [double]"27311.7998066416"
$ws.Cells.Item(2, 18).Value = [double]"245806.1982597744"
$ws.Cells.Item(2, 19).Value = [double]"0.6642269562304284"
$ws.Cells.Item(2, 20).Value = [double]"0.6642269562304284"
$ws.Cells.Item(3, 7).Value = [double]"60.813934"
$ws.Cells.Item(3, 8).Value = [double]"182.441802"
$ws.Cells.Item(3, 9).Value = [double]"0.6840634102070431"
$ws.Cells.Item(3, 10).Value = [double]"0.6840634102070431"
$ws.Cells.Item(3, 15).Value = [double]"0.01131353526791385"
$ws.Cells.Item(3, 16).Value = [double]"0.01131353526791385"
$ws.Cells.Item(3, 17).Value = [double]"318.22076837213"
$ws.Cells.Item(3, 18).Value = [double]"2863.98691534917"
$ws.Cells.Item(3, 19).Value = [double]"0.007739175516866798"
$ws.Cells.Item(3, 20).Value = [double]"0.007739175516866799"
$ws.Cells.Item(4, 7).Value = [double]"60.813934"
$ws.Cells.Item(4, 8).Value = [double]"182.441802"
$ws.Cells.Item(4, 9).Value = [double]"0.6840634102070431"
$ws.Cells.Item(4, 10).Value = [double]"0.6840634102070431"
$ws.Cells.Item(4, 13).Value = [double]"5.266527"
$ws.Cells.Item(4, 14).Value = [double]"15.799581"
$ws.Cells.Item(4, 15).Value = [double]"0.0113866829528418"
$ws.Cells.Item(4, 16).Value = [double]"0.0113866829528418"
$ws.Cells.Item(4, 17).Value = [double]"320.278225387218"
$ws.Cells.Item(4, 18).Value = [double]"2882.504028484962"
$ws.Cells.Item(4, 19).Value = [double]"0.007789213171667363"
$ws.Cells.Item(4, 20).Value = [double]"0.007789213171667363"
$ws.Cells.Item(5, 7).Value = [double]"60.813934"
$ws.Cells.Item(5, 8).Value = [double]"182.441802"
$ws.Cells.Item(5, 9).Value = [double]"0.6840634102070431"
$ws.Cells.Item(5, 10).Value = [double]"0.6840634102070431"
$ws.Cells.Item(5, 13).Value = [double]"2.912815666666667"
$ws.Cells.Item(5, 14).Value = [double]"8.738447000000001"
$ws.Cells.Item(5, 15).Value = [double]"0.006297757230980464"
$ws.Cells.Item(5, 16).Value = [double]"0.006297757230980464"
$ws.Cells.Item(5, 17).Value = [double]"177.1397797068327"
$ws.Cells.Item(5, 18).Value = [double]"1594.258017361494"
$ws.Cells.Item(5, 19).Value = [double]"0.004308065288080561"
$ws.Cells.Item(5, 20).Value = [double]"0.004308065288080561"
$ws.Cells.Item(6, 7).Value = [double]"2.823821666666666"
$ws.Cells.Item(6, 8).Value = [double]"8.471464999999998"
$ws.Cells.Item(6, 9).Value = [double]"0.0317636592810545"
$ws.Cells.Item(6, 10).Value = [double]"0.0317636592810545"
$ws.Cells.Item(6, 13).Value = [double]"449.104309"
$ws.Cells.Item(6, 14).Value = [double]"1347.312927"
$ws.Cells.Item(6, 15).Value = [double]"0.9710020245482639"
$ws.Cells.Item(6, 16).Value = [double]"0.9710020245482639"
$ws.Cells.Item(6, 17).Value = [double]"1268.190478347561"
$ws.Cells.Item(6, 18).Value = [double]"11413.71430512805"
$ws.Cells.Item(6, 19).Value = [double]"0.03084257746896517"
$ws.Cells.Item(6, 20).Value = [double]"0.03084257746896517"
$ws.Cells.Item(7, 7).Value = [double]"2.823821666666666"
$ws.Cells.Item(7, 8).Value = [double]"8.471464999999998"
$ws.Cells.Item(7, 9).Value = [double]"0.0317636592810545"
$ws.Cells.Item(7, 10).Value = [double]"0.0317636592810545"
$ws.Cells.Item(7, 15).Value = [double]"0.01131353526791385"
$ws.Cells.Item(7, 16).Value = [double]"0.01131353526791385"
$ws.Cells.Item(7, 19).Value = [double]"0.000359359279514209"
$ws.Cells.Item(7, 20).Value = [double]"0.000359359279514209"
$ws.Cells.Item(8, 7).Value = [double]"2.823821666666666"
$ws.Cells.Item(8, 8).Value = [double]"8.471464999999998"
$ws.Cells.Item(8, 9).Value = [double]"0.0317636592810545"
$ws.Cells.Item(8, 10).Value = [double]"0.0317636592810545"
$ws.Cells.Item(8, 13).Value = [double]"5.266527"
$ws.Cells.Item(8, 14).Value = [double]"15.799581"
$ws.Cells.Item(8, 15).Value = [double]"0.0113866829528418"
$ws.Cells.Item(8, 16).Value = [double]"0.0113866829528418"
$ws.Cells.Item(8, 17).Value = [double]"14.871733050685"
$ws.Cells.Item(8, 18).Value = [double]"133.845597456165"
$ws.Cells.Item(8, 19).Value = [double]"0.0003616827176554584"
$ws.Cells.Item(8, 20).Value = [double]"0.0003616827176554584"
$ws.Cells.Item(9, 7).Value = [double]"2.823821666666666"
$ws.Cells.Item(9, 8).Value = [double]"8.471464999999998"
$ws.Cells.Item(9, 9).Value = [double]"0.0317636592810545"
$ws.Cells.Item(9, 10).Value = [double]"0.0317636592810545"
$ws.Cells.Item(9, 13).Value = [double]"2.912815666666667"
$ws.Cells.Item(9, 14).Value = [double]"8.738447000000001"
$ws.Cells.Item(9, 15).Value = [double]"0.006297757230980464"
$ws.Cells.Item(9, 16).Value = [double]"0.006297757230980464"
$ws.Cells.Item(9, 17).Value = [double]"8.225271990539444"
$ws.Cells.Item(9, 18).Value = [double]"74.02744791485499"
$ws.Cells.Item(9, 19).Value = [double]"0.0002000398149196607"
$ws.Cells.Item(9, 20).Value = [double]"0.0002000398149196607"
$ws.Cells.Item(10, 5).Value = [double]"2"
$ws.Cells.Item(10, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(10, 7).Value = [double]"0.640208"
$ws.Cells.Item(10, 8).Value = [double]"1.920624"
$ws.Cells.Item(10, 9).Value = [double]"0.007201357302782462"
$ws.Cells.Item(10, 10).Value = [double]"0.007201357302782463"
$ws.Cells.Item(10, 13).Value = [double]"449.104309"
$ws.Cells.Item(10, 14).Value = [double]"1347.312927"
$ws.Cells.Item(10, 15).Value = [double]"0.9710020245482639"
$ws.Cells.Item(10, 16).Value = [double]"0.9710020245482639"
$ws.Cells.Item(10, 17).Value = [double]"287.520171456272"
$ws.Cells.Item(10, 18).Value = [double]"2587.681543106448"
$ws.Cells.Item(10, 19).Value = [double]"0.006992532520497195"
$ws.Cells.Item(10, 20).Value = [double]"0.006992532520497196"
$ws.Cells.Item(11, 5).Value = [double]"2"
$ws.Cells.Item(11, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(11, 7).Value = [double]"0.640208"
$ws.Cells.Item(11, 8).Value = [double]"1.920624"
$ws.Cells.Item(11, 9).Value = [double]"0.007201357302782462"
$ws.Cells.Item(11, 10).Value = [double]"0.007201357302782463"
$ws.Cells.Item(11, 15).Value = [double]"0.01131353526791385"
$ws.Cells.Item(11, 16).Value = [double]"0.01131353526791385"
$ws.Cells.Item(11, 17).Value = [double]"3.35001320056"
$ws.Cells.Item(11, 18).Value = [double]"30.15011880504"
$ws.Cells.Item(11, 19).Value = [double]"8.147280982187831E-05"
$ws.Cells.Item(11, 20).Value = [double]"8.147280982187832E-05"
$ws.Cells.Item(12, 5).Value = [double]"2"
$ws.Cells.Item(12, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(12, 7).Value = [double]"0.640208"
$ws.Cells.Item(12, 8).Value = [double]"1.920624"
$ws.Cells.Item(12, 9).Value = [double]"0.007201357302782462"
$ws.Cells.Item(12, 10).Value = [double]"0.007201357302782463"
$ws.Cells.Item(12, 13).Value = [double]"5.266527"
$ws.Cells.Item(12, 14).Value = [double]"15.799581"
$ws.Cells.Item(12, 15).Value = [double]"0.0113866829528418"
$ws.Cells.Item(12, 16).Value = [double]"0.0113866829528418"
$ws.Cells.Item(12, 17).Value = [double]"3.371672717616"
$ws.Cells.Item(12, 18).Value = [double]"30.345054458544"
$ws.Cells.Item(12, 19).Value = [double]"8.199957243691585E-05"
$ws.Cells.Item(12, 20).Value = [double]"8.199957243691585E-05"
$ws.Cells.Item(13, 5).Value = [double]"2"
$ws.Cells.Item(13, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(13, 7).Value = [double]"0.640208"
$ws.Cells.Item(13, 8).Value = [double]"1.920624"
$ws.Cells.Item(13, 9).Value = [double]"0.007201357302782462"
$ws.Cells.Item(13, 10).Value = [double]"0.007201357302782463"
$ws.Cells.Item(13, 13).Value = [double]"2.912815666666667"
$ws.Cells.Item(13, 14).Value = [double]"8.738447000000001"
$ws.Cells.Item(13, 15).Value = [double]"0.006297757230980464"
$ws.Cells.Item(13, 16).Value = [double]"0.006297757230980464"
$ws.Cells.Item(13, 17).Value = [double]"1.864807892325334"
$ws.Cells.Item(13, 18).Value = [double]"16.783271030928"
$ws.Cells.Item(13, 19).Value = [double]"4.535240002647222E-05"
$ws.Cells.Item(13, 20).Value = [double]"4.535240002647222E-05"
$ws.Cells.Item(14, 7).Value = [double]"24.623055"
$ws.Cells.Item(14, 8).Value = [double]"73.86916500000001"
$ws.Cells.Item(14, 9).Value = [double]"0.2769715732091199"
$ws.Cells.Item(14, 10).Value = [double]"0.2769715732091199"
$ws.Cells.Item(14, 13).Value = [double]"449.104309"
$ws.Cells.Item(14, 14).Value = [double]"1347.312927"
$ws.Cells.Item(14, 15).Value = [double]"0.9710020245482639"
$ws.Cells.Item(14, 16).Value = [double]"0.9710020245482639"
$ws.Cells.Item(14, 17).Value = [double]"11058.320101244"
$ws.Cells.Item(14, 18).Value = [double]"99524.88091119597"
$ws.Cells.Item(14, 19).Value = [double]"0.2689399583283731"
$ws.Cells.Item(14, 20).Value = [double]"0.2689399583283731"
$ws.Cells.Item(15, 7).Value = [double]"24.623055"
$ws.Cells.Item(15, 8).Value = [double]"73.86916500000001"
$ws.Cells.Item(15, 9).Value = [double]"0.2769715732091199"
$ws.Cells.Item(15, 10).Value = [double]"0.2769715732091199"
$ws.Cells.Item(15, 15).Value = [double]"0.01131353526791385"
$ws.Cells.Item(15, 16).Value = [double]"0.01131353526791385"
$ws.Cells.Item(15, 17).Value = [double]"128.844936783225"
$ws.Cells.Item(15, 18).Value = [double]"1159.604431049025"
$ws.Cells.Item(15, 19).Value = [double]"0.00313352766171096"
$ws.Cells.Item(15, 20).Value = [double]"0.00313352766171096"
$ws.Cells.Item(16, 7).Value = [double]"24.623055"
$ws.Cells.Item(16, 8).Value = [double]"73.86916500000001"
$ws.Cells.Item(16, 9).Value = [double]"0.2769715732091199"
$ws.Cells.Item(16, 10).Value = [double]"0.2769715732091199"
$ws.Cells.Item(16, 13).Value = [double]"5.266527"
$ws.Cells.Item(16, 14).Value = [double]"15.799581"
$ws.Cells.Item(16, 15).Value = [double]"0.0113866829528418"
$ws.Cells.Item(16, 16).Value = [double]"0.0113866829528418"
$ws.Cells.Item(16, 17).Value = [double]"129.677983979985"
$ws.Cells.Item(16, 18).Value = [double]"1167.101855819865"
$ws.Cells.Item(16, 19).Value = [double]"0.003153787491082059"
$ws.Cells.Item(16, 20).Value = [double]"0.003153787491082059"
$ws.Cells.Item(17, 7).Value = [double]"24.623055"
$ws.Cells.Item(17, 8).Value = [double]"73.86916500000001"
$ws.Cells.Item(17, 9).Value = [double]"0.2769715732091199"
$ws.Cells.Item(17, 10).Value = [double]"0.2769715732091199"
$ws.Cells.Item(17, 13).Value = [double]"2.912815666666667"
$ws.Cells.Item(17, 14).Value = [double]"8.738447000000001"
$ws.Cells.Item(17, 15).Value = [double]"0.006297757230980464"
$ws.Cells.Item(17, 16).Value = [double]"0.006297757230980464"
$ws.Cells.Item(17, 17).Value = [double]"71.72242036519502"
$ws.Cells.Item(17, 18).Value = [double]"645.5017832867551"
$ws.Cells.Item(17, 19).Value = [double]"0.00174429972795377"
$ws.Cells.Item(17, 20).Value = [double]"0.00174429972795377"
